$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75, shifting existing rows 75-84 down to 76-85
$ws.Rows.Item(75).Insert()

$ws.Cells.Item(75, 1).Value = 3
$ws.Cells.Item(75, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(75, 3).Value = "Coquimbo"
$ws.Cells.Item(75, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(75, 4).Value = 45223
$ws.Cells.Item(75, 5).Value = 5
$ws.Cells.Item(75, 6).Value = 300000000
$ws.Cells.Item(75, 7).Value = "Espárragos"
$ws.Cells.Item(75, 8).Value = "Verde"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 1500
$ws.Cells.Item(75, 11).Value = 1700
$ws.Cells.Item(75, 12).Value = 1700
$ws.Cells.Item(75, 13).Value = 1700
$ws.Cells.Item(75, 14).Value = "$/kilo"
$ws.Cells.Item(75, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(75, 16).Value = 1700
$ws.Cells.Item(75, 17).Value = 1
$ws.Cells.Item(75, 18).Value = "Hortaliza"
